$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 8492.789000000001
$ws.Range("I76").Value = 8736.200000000001
$ws.Range("K76").Value = 8736.200000000001
$ws.Range("M76").Value = -8421.200000000001
$ws.Range("H79").Value = 8492.789000000001
$ws.Range("I79").Value = 8736.200000000001
$ws.Range("K79").Value = 8736.200000000001
$ws.Range("M79").Value = -7644.200000000001
$ws.Range("H98").Value = 188757
$ws.Range("J98").Value = 1503323.8
$ws.Range("L98").Value = 1503323.8
$ws.Range("N98").Value = -1506319.8
$ws.Range("H122").Value = 188757
$ws.Range("J122").Value = 1503323.8
$ws.Range("L122").Value = 4509971.4
$ws.Range("N122").Value = -4514871.4
$ws.Range("H138").Value = 5304.5293
$ws.Range("J138").Value = 5401.1377
$ws.Range("L138").Value = 16203.4131
$ws.Range("N138").Value = -26483.4131
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 5951.636
$ws.Range("I2").Value = 1551.8462
$ws.Range("J2").Value = 12306.889
$ws.Range("K2").Value = 1551.8462
$ws.Range("L2").Value = 12306.889
$ws.Range("M2").Value = -1438.8462
$ws.Range("N2").Value = -12532.889
$ws.Range("H32").Value = 7350.1313
$ws.Range("I32").Value = 6655.7334
$ws.Range("K32").Value = 6655.7334
$ws.Range("M32").Value = -6368.7334
$ws.Range("H61").Value = 4532.5747
$ws.Range("I61").Value = 4139.683
$ws.Range("K61").Value = 4139.683
$ws.Range("M61").Value = -3927.683
$ws.Range("H74").Value = 3208.739
$ws.Range("I74").Value = 2116.9412
$ws.Range("K74").Value = 2116.9412
$ws.Range("M74").Value = -1242.9412
$ws.Range("H77").Value = 3208.739
$ws.Range("I77").Value = 2116.9412
$ws.Range("K77").Value = 10584.706
$ws.Range("M77").Value = -6216.706000000002
$ws.Range("H116").Value = 5951.636
$ws.Range("I116").Value = 1551.8462
$ws.Range("J116").Value = 12306.889
$ws.Range("K116").Value = 1551.8462
$ws.Range("L116").Value = 12306.889
$ws.Range("M116").Value = 742.1538
$ws.Range("N116").Value = -16894.889
$ws.Range("H136").Value = 4532.5747
$ws.Range("I136").Value = 4139.683
$ws.Range("K136").Value = 12419.049
$ws.Range("M136").Value = -9869.048999999999
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 5951.636
$ws.Range("I3").Value = 1551.8462
$ws.Range("J3").Value = 12306.889
$ws.Range("K3").Value = 1551.8462
$ws.Range("L3").Value = 12306.889
$ws.Range("M3").Value = -1437.8462
$ws.Range("N3").Value = -12534.889
$ws.Range("H20").Value = 4397.6665
$ws.Range("I20").Value = 3592.3333
$ws.Range("J20").Value = 5203
$ws.Range("K20").Value = 3592.3333
$ws.Range("L20").Value = 5203
$ws.Range("M20").Value = -3345.3333
$ws.Range("N20").Value = -5697
$ws.Range("H105").Value = 21572.234
$ws.Range("I105").Value = 23105.1
$ws.Range("K105").Value = 23105.1
$ws.Range("M105").Value = -21358.1
$ws.Range("H117").Value = 73682
$ws.Range("J117").Value = 73682
$ws.Range("L117").Value = 73682
$ws.Range("N117").Value = -82860
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 2657.0303
$ws.Range("I107").Value = 2550.0454
$ws.Range("J107").Value = 2871
$ws.Range("K107").Value = 2550.0454
$ws.Range("L107").Value = 2871
$ws.Range("M107").Value = -630.0454
$ws.Range("N107").Value = -6711
$ws.Range("H122").Value = 6636.6
$ws.Range("I122").Value = 2570.3333
$ws.Range("J122").Value = 9347.444
$ws.Range("K122").Value = 7710.999899999999
$ws.Range("L122").Value = 28042.332
$ws.Range("M122").Value = -5260.999899999999
$ws.Range("N122").Value = -32942.33199999999
$ws.Range("H134").Value = 2824.1853
$ws.Range("I134").Value = 1884.9584
$ws.Range("K134").Value = 5654.8752
$ws.Range("M134").Value = -3119.8752
$ws.Range("H139").Value = 100154.5
$ws.Range("J139").Value = 149600
$ws.Range("L139").Value = 149600
$ws.Range("N139").Value = -159880
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1853084.2
$ws.Range("J5").Value = 4167827.5
$ws.Range("L5").Value = 12503482.5
$ws.Range("N5").Value = -12503706.5
$ws.Range("H22").Value = 336633.66
$ws.Range("I22").Value = 899
$ws.Range("J22").Value = 504501
$ws.Range("K22").Value = 2697
$ws.Range("L22").Value = 1513503
$ws.Range("M22").Value = -2528
$ws.Range("N22").Value = -1513841
$ws.Range("H27").Value = 336633.66
$ws.Range("I27").Value = 899
$ws.Range("J27").Value = 504501
$ws.Range("K27").Value = 2697
$ws.Range("L27").Value = 1513503
$ws.Range("M27").Value = -2595
$ws.Range("N27").Value = -1513707
$ws.Range("H44").Value = 125119
$ws.Range("I44").Value = 100
$ws.Range("K44").Value = 300
$ws.Range("M44").Value = 98
$ws.Range("H68").Value = 3574.8333
$ws.Range("I68").Value = 326
$ws.Range("J68").Value = 4224.6
$ws.Range("K68").Value = 978
$ws.Range("L68").Value = 12673.8
$ws.Range("M68").Value = -167
$ws.Range("N68").Value = -14295.8
$ws.Range("H71").Value = 3574.8333
$ws.Range("I71").Value = 326
$ws.Range("J71").Value = 4224.6
$ws.Range("K71").Value = 2934
$ws.Range("L71").Value = 38021.4
$ws.Range("M71").Value = 1122
$ws.Range("N71").Value = -46133.4
$ws.Range("H98").Value = 3501
$ws.Range("J98").Value = 3501
$ws.Range("L98").Value = 10503
$ws.Range("N98").Value = -13499
$ws.Range("H113").Value = 2648.375
$ws.Range("I113").Value = 2396.8
$ws.Range("J113").Value = 2762.7273
$ws.Range("K113").Value = 7190.400000000001
$ws.Range("L113").Value = 8288.1819
$ws.Range("M113").Value = -5020.400000000001
$ws.Range("N113").Value = -12628.1819
$ws.Range("H122").Value = 2151.1667
$ws.Range("I122").Value = 655.8333
$ws.Range("J122").Value = 2898.8333
$ws.Range("K122").Value = 5902.4997
$ws.Range("L122").Value = 26089.4997
$ws.Range("M122").Value = -3452.4997
$ws.Range("N122").Value = -30989.4997
$ws.Range("H131").Value = 5325939.5
$ws.Range("J131").Value = 3413233.5
$ws.Range("L131").Value = 10239700.5
$ws.Range("N131").Value = -10249780.5
$ws.Range("H132").Value = 4586.56
$ws.Range("J132").Value = 5033.7856
$ws.Range("L132").Value = 45304.0704
$ws.Range("N132").Value = -50364.0704
$ws.Range("H135").Value = 1853084.2
$ws.Range("J135").Value = 4167827.5
$ws.Range("L135").Value = 37510447.5
$ws.Range("N135").Value = -37515517.5
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 361639.44
$ws.Range("I80").Value = 558216.5600000001
$ws.Range("J80").Value = 7800.6
$ws.Range("K80").Value = 558216.5600000001
$ws.Range("L80").Value = 7800.6
$ws.Range("M80").Value = -557218.5600000001
$ws.Range("N80").Value = -9796.6
$ws.Range("H83").Value = 361639.44
$ws.Range("I83").Value = 558216.5600000001
$ws.Range("J83").Value = 7800.6
$ws.Range("K83").Value = 2791082.8
$ws.Range("L83").Value = 39003
$ws.Range("M83").Value = -2786090.8
$ws.Range("N83").Value = -48987
$ws.Range("H122").Value = 6328.125
$ws.Range("I122").Value = 4489.636
$ws.Range("K122").Value = 13468.908
$ws.Range("M122").Value = -11018.908
$ws.Range("H132").Value = 2255.04
$ws.Range("I132").Value = 1556.75
$ws.Range("K132").Value = 4670.25
$ws.Range("M132").Value = -2140.25
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 2862.7368
$ws.Range("I46").Value = 1448
$ws.Range("K46").Value = 1448
$ws.Range("M46").Value = -1260
$ws.Range("H68").Value = 7049.5
$ws.Range("I68").Value = 3873.5
$ws.Range("K68").Value = 3873.5
$ws.Range("M68").Value = -3124.5
$ws.Range("H71").Value = 7049.5
$ws.Range("I71").Value = 3873.5
$ws.Range("K71").Value = 19367.5
$ws.Range("M71").Value = -15623.5
$ws.Range("H93").Value = 5332.3335
$ws.Range("I93").Value = 4001.5
$ws.Range("K93").Value = 4001.5
$ws.Range("M93").Value = -2753.5
$ws.Range("H136").Value = 5496.0713
$ws.Range("I136").Value = 3317.5
$ws.Range("K136").Value = 9952.5
$ws.Range("M136").Value = -7402.5
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 986.8293
$ws.Range("I107").Value = 815.56665
$ws.Range("J107").Value = 1453.909
$ws.Range("K107").Value = 2446.69995
$ws.Range("L107").Value = 4361.727000000001
$ws.Range("M107").Value = -526.6999500000002
$ws.Range("N107").Value = -8201.727000000001
$ws.Range("H116").Value = 183223.5
$ws.Range("J116").Value = 183223.5
$ws.Range("L116").Value = 183223.5
$ws.Range("N116").Value = -192401.5
$ws.Range("H122").Value = 4451
$ws.Range("I122").Value = 3101.8
$ws.Range("K122").Value = 9305.400000000001
$ws.Range("M122").Value = -6855.400000000001
$ws.Range("H136").Value = 2792.0657
$ws.Range("I136").Value = 1506.3478
$ws.Range("J136").Value = 6734.933
$ws.Range("K136").Value = 4519.0434
$ws.Range("L136").Value = 20204.799
$ws.Range("M136").Value = -1969.0434
$ws.Range("N136").Value = -25304.799
